$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A106").Value = "B0B7VF1S57"
$ws.Range("B106").Value = "WN-KXD6-HSUQ"
$ws.Range("E106").Value = "Sognare® Almohada Fussión 6 Pack Tamaño Estándar, Relleno Suave, Anti Ácaros y Lavable. Hipoalergénico, Anti Ácaros y Lavable. Cont. 6 Pzas."
$ws.Range("F106").Value = "Almohada"

$ws.Range("A107").Value = "B0DFKVST8R"
$ws.Range("B107").Value = "UU-JW8I-VAU7"
$ws.Range("E107").Value = "Sognare Theramed PRO Colchoneta de Masaje con 3 Terapias: Calor Infrarrojo, Vibración y Masaje Shiatsu. Alivia la Tensión, Rigidez y Dolor Corporal. 100 Noches de Garantía."
$ws.Range("F107").Value = "Theramed"

$ws.Range("A108").Value = "B0B76D5HF2"
$ws.Range("B108").Value = "6V-J7S8-F5H6"
$ws.Range("E108").Value = "Sognare Set 1 Cubre Colchón Individual + 1 Almohada Estandar Fussión Firme + 1 Almohada Suave. Hipoalergénico, Anti Ácaros y Lavable. Cont. 3 Pzas."
$ws.Range("F108").Value = "Cubre"

$ws.Range("A109").Value = "B0B7QN1K82"
$ws.Range("B109").Value = "8B-9WVK-ISVT"
$ws.Range("E109").Value = "Sognare Set 1 Cubre Colchón King Size Extra Confort + 2 Almohadas Estándar Fussion Relleno Suave, 100% Algodón, Termorregulable, Hipoalergenico, Anti ácaros. Cont. 3 pzas."
$ws.Range("F109").Value = "Cubre"

# Match the saved view state: frozen pane scrolled near the new rows, with
# B107 as the active/selected cell.
$win = $excel.ActiveWindow
$win.FreezePanes = $true
[void]$ws.Range("B107").Select()
